# Fill in the team names and roster members for CS320 Sp24 Teams sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Section 101 table ---------------------------------------------------
# Header row (team names / project titles)
$ws.Range("B6").Value = "Team 101-1:`nRevMetrix Bowler UI"
$ws.Range("C6").Value = "Team 101-2:`nProductivity Planner"
$ws.Range("D6").Value = "Team 101-3:`nTBAG: York County Ghosts"
$ws.Range("E6").Value = "Team 101-4:`n<TBD>"

# Team members
$ws.Range("B7").Value = "Brandon Woodward"
$ws.Range("C7").Value = "Deborah Amao"
$ws.Range("D7").Value = "Joshua Byers"

$ws.Range("B8").Value = "Emily Culp"
$ws.Range("C8").Value = "Kevin Lindemann"
$ws.Range("D8").Value = "Ren De Alva"

$ws.Range("B9").Value = "Emmet Larson"
$ws.Range("C9").Value = "Ryan Huber"
$ws.Range("D9").Value = "Thomas Wakeland"

$ws.Range("B10").Value = "Zachary Cox"
$ws.Range("D10").Value = "Ethan VonStein"

# --- Section 102 table ---------------------------------------------------
# Header row (team names / project titles)
$ws.Range("B14").Value = "Team 102-1:`nTBAG: Tea-Bag"
$ws.Range("C14").Value = "Team 102-2:`nRevMetrix Bowler UI"
$ws.Range("D14").Value = "Team 102-3:`nTBAG: Apocrypha"
$ws.Range("E14").Value = "Team 102-4:`n<TBD>"

# Team members
$ws.Range("B15").Value = "Andrew Loiseau"
$ws.Range("C15").Value = "Charles Carroll"
$ws.Range("D15").Value = "Carson Mack"

$ws.Range("B16").Value = "Alyssa Nelson"
$ws.Range("C16").Value = "Gabriel Manero"
$ws.Range("D16").Value = "Jonathan Waight"

$ws.Range("B17").Value = "Matthew Brown"
$ws.Range("C17").Value = "Jakeb Nielsen"
$ws.Range("D17").Value = "Korbin Dick"

$ws.Range("C18").Value = "Ryon Washington"
$ws.Range("D18").Value = "Spencer Hayes"

# --- View state -----------------------------------------------------------
$ws.Range("D6").Select()
